# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.865.07"
$ws.Range("E2").Value = "  +3.89%  "

$ws.Range("D3").Value = "1.864.16"
$ws.Range("E3").Value = "  +2.76%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.33"
$ws.Range("E5").Value = "  +2.29%  "

$ws.Range("E6").Value = "  +2.95%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.68"
$ws.Range("E8").Value = "  +12.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.312"
$ws.Range("E9").Value = "  +7.16%  "

$ws.Range("E10").Value = "  +3.18%  "

$ws.Range("E11").Value = "  +3.43%  "

$ws.Range("D12").Value = "2.130.23"
$ws.Range("E12").Value = "  +2.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.58"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("D14").Value = "1.863.27"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.681"
$ws.Range("E15").Value = "  +7.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.80"
$ws.Range("E16").Value = "  +8.13%  "

$ws.Range("D17").Value = "35.789.29"
$ws.Range("E17").Value = "  +3.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.37"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("D19").Value = "0.0₃0805"
$ws.Range("E19").Value = "  +3.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.98"
$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.32"
$ws.Range("E21").Value = "  +8.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.81"
$ws.Range("E22").Value = "  +16.04%  "

$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.48"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.05"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.96"
$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("E28").Value = "  +2.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("E29").Value = "  +17.91%  "

$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "3.338.89"
$ws.Range("E31").Value = "  +37.42%  "

$ws.Range("E32").Value = "  +3.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0546"
$ws.Range("E33").Value = "  +4.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.07"
$ws.Range("E34").Value = "  +5.45%  "

$ws.Range("E35").Value = "  +3.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "97.56"
$ws.Range("E36").Value = "  +19.24%  "

$ws.Range("E37").Value = "  +6.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.53"
$ws.Range("E38").Value = "  +7.57%  "

$ws.Range("E39").Value = "  +2.40%  "

$ws.Range("D40").Value = "1.346.96"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0197"
$ws.Range("E41").Value = "  +4.94%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.03"
$ws.Range("E42").Value = "  +7.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.21"
$ws.Range("E43").Value = "  +10.27%  "

$ws.Range("E44").Value = "  +5.22%  "

$ws.Range("E45").Value = "  +2.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("E47").Value = "  +2.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.12"
$ws.Range("E48").Value = "  +5.52%  "

$ws.Range("D49").Value = "2.028.43"
$ws.Range("E49").Value = "  +2.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.99"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("E51").Value = "  +0.19%  "
